# Append a new row (row 20) of forecast data to the active worksheet,
# mirroring the style/format of the existing rows above it (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 20
$prevRow = $newRow - 1

# Duplicate the formatting (style/number format) of the cell directly above
# the new date cell (A19 -> A20) so the new row looks consistent with the rest
# of the table.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

# A20: date serial value (2025-11-25)
$ws.Cells.Item($newRow, 1).Value = 45986

# B20: year value
$ws.Cells.Item($newRow, 2).Value = 2025

# C20: forecast value
$ws.Cells.Item($newRow, 3).Value = -0.08656168856399082

# D20: year value
$ws.Cells.Item($newRow, 4).Value = 2026

# E20: forecast value
$ws.Cells.Item($newRow, 5).Value = 0.6232357314897463
